# Adding few test cases along with new custom utilities
# - keep the existing "loginData" sheet, but move the current selection
# - add a new "pimAddEmpoyeeData" worksheet after it with employee name test data

$wb = $excel.ActiveWorkbook

# The user had re-selected a different cell on loginData before switching
# to the newly added sheet, so loginData ends up with selection F11 and
# loses "tabSelected".
$ws1 = $wb.Worksheets.Item("loginData")
$ws1.Activate() | Out-Null
$ws1.Range("F11").Select() | Out-Null

# Add the new worksheet right after the last existing sheet (loginData)
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "pimAddEmpoyeeData"

# Header row
$ws2.Range("A1").Value = "FirstName"
$ws2.Range("B1").Value = "MiddleName"
$ws2.Range("C1").Value = "LastName"

# Test data rows
$ws2.Range("A2").Value = "Srinivasa"
$ws2.Range("B2").Value = "N"
$ws2.Range("C2").Value = "Raj"

$ws2.Range("A3").Value = "Charan"
$ws2.Range("B3").Value = "B"
$ws2.Range("C3").Value = "Kumar"

# Widen the first two columns so the names are fully visible
$ws2.Columns.Item(1).ColumnWidth = 17.125
$ws2.Columns.Item(2).ColumnWidth = 17.43

# Leave the cursor where the author left it on the new sheet
$ws2.Range("E19").Select() | Out-Null
